# Car_type_Mapping / mapping details.xlsx
# Adds a new "PALM OIL" / "OILS,NUT,SEED" mapping row beneath the existing
# Car Name / Commodity Name table, matching the formatting already used by
# the rest of the data rows (copied from row 11), and leaves the selection
# on D14 as it was when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: column A = Car Name, column B = Commodity Name.
$ws.Range("B13").Value = "OILS,NUT,SEED"
$ws.Range("A13").Value = "PALM OIL"

# Match the look of the other data rows (border + wrap text) by copying the
# formatting from the row above (A11:B11 already carries that style).
$ws.Range("A11:B11").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author left it.
$ws.Range("D14").Select() | Out-Null
